# chore: update Sheets via scheduled runner
# Refreshes market-price / profit figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leve rows across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 11780.417
$ws.Range("I74").Value = 11820.571
$ws.Range("K74").Value = 11820.571
$ws.Range("M74").Value = -10884.571
$ws.Range("H76").Value = 4334.3335
$ws.Range("J76").Value = 4004
$ws.Range("L76").Value = 4004
$ws.Range("N76").Value = -4634
$ws.Range("H77").Value = 11780.417
$ws.Range("I77").Value = 11820.571
$ws.Range("K77").Value = 59102.855
$ws.Range("M77").Value = -54422.855
$ws.Range("H79").Value = 4334.3335
$ws.Range("J79").Value = 4004
$ws.Range("L79").Value = 4004
$ws.Range("N79").Value = -6188
$ws.Range("H116").Value = 3103.6
$ws.Range("I116").Value = 2962.3845
$ws.Range("K116").Value = 2962.3845
$ws.Range("M116").Value = 479.6154999999999
$ws.Range("H137").Value = 3147.9167
$ws.Range("J137").Value = 2961.125
$ws.Range("L137").Value = 8883.375
$ws.Range("N137").Value = -13983.375
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4704.7856
$ws.Range("I45").Value = 3985.3333
$ws.Range("K45").Value = 3985.3333
$ws.Range("M45").Value = -3608.3333
$ws.Range("H50").Value = 8774
$ws.Range("J50").Value = 8774
$ws.Range("L50").Value = 8774
$ws.Range("N50").Value = -10202
$ws.Range("H88").Value = 2111.625
$ws.Range("J88").Value = 2286
$ws.Range("L88").Value = 2286
$ws.Range("N88").Value = -3098
$ws.Range("H91").Value = 2111.625
$ws.Range("J91").Value = 2286
$ws.Range("L91").Value = 2286
$ws.Range("N91").Value = -5094
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3070.353
$ws.Range("I20").Value = 3156.2144
$ws.Range("K20").Value = 3156.2144
$ws.Range("M20").Value = -2909.2144
$ws.Range("H86").Value = 2782849.5
$ws.Range("I86").Value = 6064144
$ws.Range("J86").Value = 6369.769
$ws.Range("K86").Value = 6064144
$ws.Range("L86").Value = 6369.769
$ws.Range("M86").Value = -6063021
$ws.Range("N86").Value = -8615.769
$ws.Range("H89").Value = 2782849.5
$ws.Range("I89").Value = 6064144
$ws.Range("J89").Value = 6369.769
$ws.Range("K89").Value = 30320720
$ws.Range("L89").Value = 31848.845
$ws.Range("M89").Value = -30315104
$ws.Range("N89").Value = -43080.845
$ws.Range("H94").Value = 3332.926
$ws.Range("I94").Value = 980.9091
$ws.Range("K94").Value = 980.9091
$ws.Range("M94").Value = -529.9091
$ws.Range("H105").Value = 5605.75
$ws.Range("I105").Value = 3610
$ws.Range("K105").Value = 3610
$ws.Range("M105").Value = -1863
$ws.Range("H124").Value = 79495
$ws.Range("J124").Value = 79495
$ws.Range("L124").Value = 79495
$ws.Range("N124").Value = -89315
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2417.9167
$ws.Range("I16").Value = 1732.6666
$ws.Range("J16").Value = 2646.3333
$ws.Range("K16").Value = 1732.6666
$ws.Range("L16").Value = 2646.3333
$ws.Range("M16").Value = -1445.6666
$ws.Range("N16").Value = -3220.3333
$ws.Range("H31").Value = 2135.0667
$ws.Range("I31").Value = 2135.0667
$ws.Range("K31").Value = 2135.0667
$ws.Range("M31").Value = -1840.0667
$ws.Range("H34").Value = 2135.0667
$ws.Range("I34").Value = 2135.0667
$ws.Range("K34").Value = 2135.0667
$ws.Range("M34").Value = -1933.0667
$ws.Range("H62").Value = 4540.1
$ws.Range("I62").Value = 4267.3335
$ws.Range("K62").Value = 4267.3335
$ws.Range("M62").Value = -3643.3335
$ws.Range("H65").Value = 4540.1
$ws.Range("I65").Value = 4267.3335
$ws.Range("K65").Value = 21336.6675
$ws.Range("M65").Value = -18216.6675
$ws.Range("H113").Value = 2417.9167
$ws.Range("I113").Value = 1732.6666
$ws.Range("J113").Value = 2646.3333
$ws.Range("K113").Value = 1732.6666
$ws.Range("L113").Value = 2646.3333
$ws.Range("M113").Value = 437.3334
$ws.Range("N113").Value = -6986.3333
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 76.04000000000001
$ws.Range("J2").Value = 140.57143
$ws.Range("L2").Value = 843.42858
$ws.Range("N2").Value = -1069.42858
$ws.Range("H131").Value = 1583.6333
$ws.Range("J131").Value = 1600.7693
$ws.Range("L131").Value = 4802.3079
$ws.Range("N131").Value = -14882.3079
$ws.Range("H139").Value = 5504020.5
$ws.Range("I139").Value = 1113594.2
$ws.Range("J139").Value = 9096188
$ws.Range("K139").Value = 3340782.6
$ws.Range("L139").Value = 27288564
$ws.Range("M139").Value = -3335642.6
$ws.Range("N139").Value = -27298844
$ws.Range("H140").Value = 40030.133
$ws.Range("J140").Value = 4174.727
$ws.Range("L140").Value = 12524.181
$ws.Range("N140").Value = -22884.181
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 57636
$ws.Range("I57").Value = 17000
$ws.Range("J57").Value = 85573.25
$ws.Range("K57").Value = 17000
$ws.Range("L57").Value = 85573.25
$ws.Range("M57").Value = -16180
$ws.Range("N57").Value = -87213.25
$ws.Range("H80").Value = 3636.5
$ws.Range("I80").Value = 2998
$ws.Range("K80").Value = 2998
$ws.Range("M80").Value = -2000
$ws.Range("H83").Value = 3636.5
$ws.Range("I83").Value = 2998
$ws.Range("K83").Value = 14990
$ws.Range("M83").Value = -9998
$ws.Range("H113").Value = 2313.4
$ws.Range("I113").Value = 2180.7778
$ws.Range("K113").Value = 2180.7778
$ws.Range("M113").Value = -10.77779999999984
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6049.2
$ws.Range("I7").Value = 2749
$ws.Range("J7").Value = 6874.25
$ws.Range("K7").Value = 2749
$ws.Range("L7").Value = 6874.25
$ws.Range("M7").Value = -2637
$ws.Range("N7").Value = -7098.25
$ws.Range("H68").Value = 1874.1
$ws.Range("I68").Value = 1913
$ws.Range("J68").Value = 1783.3334
$ws.Range("K68").Value = 1913
$ws.Range("L68").Value = 1783.3334
$ws.Range("M68").Value = -1164
$ws.Range("N68").Value = -3281.3334
$ws.Range("H71").Value = 1874.1
$ws.Range("I71").Value = 1913
$ws.Range("J71").Value = 1783.3334
$ws.Range("K71").Value = 9565
$ws.Range("L71").Value = 8916.666999999999
$ws.Range("M71").Value = -5821
$ws.Range("N71").Value = -16404.667
$ws.Range("H126").Value = 6049.2
$ws.Range("I126").Value = 2749
$ws.Range("J126").Value = 6874.25
$ws.Range("K126").Value = 8247
$ws.Range("L126").Value = 20622.75
$ws.Range("M126").Value = -5777
$ws.Range("N126").Value = -25562.75
$ws.Range("H132").Value = 2948.194
$ws.Range("I132").Value = 2913.5
$ws.Range("J132").Value = 3035.842
$ws.Range("K132").Value = 8740.5
$ws.Range("L132").Value = 9107.526
$ws.Range("M132").Value = -6210.5
$ws.Range("N132").Value = -14167.526
$ws.Range("H140").Value = 178436.88
$ws.Range("J140").Value = 178436.88
$ws.Range("L140").Value = 178436.88
$ws.Range("N140").Value = -188796.88
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 83333.336
$ws.Range("J64").Value = 80000
$ws.Range("L64").Value = 80000
$ws.Range("N64").Value = -80496
$ws.Range("H67").Value = 83333.336
$ws.Range("J67").Value = 80000
$ws.Range("L67").Value = 80000
$ws.Range("N67").Value = -81716
$ws.Range("H96").Value = 3524.05
$ws.Range("I96").Value = 2835
$ws.Range("K96").Value = 2835
$ws.Range("M96").Value = -1462
$ws.Range("H100").Value = 828.7241
$ws.Range("I100").Value = 725.5769
$ws.Range("J100").Value = 1722.6666
$ws.Range("K100").Value = 1451.1538
$ws.Range("L100").Value = 3445.3332
$ws.Range("M100").Value = -910.1538
$ws.Range("N100").Value = -4527.3332
